$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "target" column values from "proton" to "p" for all data rows (I2:I10)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "p"
}

# Make the header row (A1:K1) bold and centered
$ws.Range("A1:K1").Font.Bold = $true
$ws.Range("A1:K1").HorizontalAlignment = -4108

# Update the selected cell to H15
$ws.Range("H15").Select()
